# Watering System Calc: swap old "ladder resistor" sheets for a low-power
# battery-budget sheet ("Power") and keep the ruler-sensor calibration
# sheet (renamed "Ruler_Sensor"). The old "Tabelle1" scratch sheet is
# removed entirely.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- 1. Drop the obsolete "Tabelle1" sheet -------------------------------
$wb.Worksheets.Item("Tabelle1").Delete() | Out-Null

# --- 2. Re-purpose "Sheet1" as "Power" -----------------------------------
$power = $wb.Worksheets.Item("Sheet1")
$power.Name = "Power"
$power.Cells.Clear() | Out-Null

$power.Cells.Item(1, 1).Value = "Battery_cap [Ah]"
$power.Cells.Item(1, 2).Formula = "=550/1000"

$power.Cells.Item(2, 1).Value = "Battery_V [V]"
$power.Cells.Item(2, 2).Value = 9

$power.Cells.Item(3, 1).Value = "Battery_P_max [Wh]"
$power.Cells.Item(3, 2).Formula = "=B2*B1"

$power.Cells.Item(5, 1).Value = "I_normal[A]"
$power.Cells.Item(5, 2).Value = 0.43

$power.Cells.Item(6, 1).Value = "I_idle[A]"
$power.Cells.Item(6, 2).Value = 0.00023

$power.Cells.Item(7, 1).Value = "normal_on_time [m]"
$power.Cells.Item(7, 2).Value = 20

$power.Cells.Item(8, 1).Value = "idle_time [m]"
$power.Cells.Item(8, 2).Formula = "=6*60"

$power.Cells.Item(9, 1).Value = "I_avg [A]"
$power.Cells.Item(9, 2).Formula = "=(B5*B7 + B6*B8)/(B7+B8)"

$power.Cells.Item(10, 1).Value = "V_supply [V]"
$power.Cells.Item(10, 2).Value = 5

$power.Cells.Item(11, 1).Value = "P_avg[W]"
$power.Cells.Item(11, 2).Formula = "=B9*B10"

$power.Cells.Item(13, 1).Value = "SMPS_eff[%]"
$power.Cells.Item(13, 2).Value = 0.9

$power.Cells.Item(14, 1).Value = "P_in[W]"
$power.Cells.Item(14, 2).Formula = "=B11/B13"

$power.Cells.Item(15, 1).Value = "Bat_lifetime [h]"
$power.Cells.Item(15, 2).Formula = "=B3/B14"

# cells that merely carry the shared scientific-notation style, no value
$power.Cells.Item(1, 5).NumberFormat = "0.00E+00"
$power.Cells.Item(4, 2).NumberFormat = "0.00E+00"
$power.Cells.Item(5, 4).NumberFormat = "0.00E+00"
$power.Cells.Item(6, 4).NumberFormat = "0.00E+00"

$power.Range("B1:B11").NumberFormat = "0.00E+00"
$power.Cells.Item(14, 2).NumberFormat = "0.00E+00"
$power.Cells.Item(15, 2).NumberFormat = "0.00E+00"

$power.Columns.Item(1).ColumnWidth = 32.71
$power.Columns.Item(2).ColumnWidth = 23.83

$power.Range("L19").Select() | Out-Null
$power.Activate() | Out-Null

# --- 3. Rename "Sheet2" (ruler/sensor calibration data) ------------------
$ruler = $wb.Worksheets.Item("Sheet2")
$ruler.Name = "Ruler_Sensor"

# --- 4. Broken external-link chart series --------------------------------
$co = $ruler.ChartObjects().Item(1)
$ser = $co.Chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(,#REF!,#REF!,1)"

$wb.Save()
